$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# This change swaps the data rows for the two files
# "8df1a12c-6253-4bfa-bf3e-77a8bf70f174.md" (row 5) and
# "969e8835-09f2-4e4a-b193-79522ee0a29a.md" (row 6) across all three
# sheets (Overview, zh-cn, de-de), and updates the status of the
# "969e8835" item (now row 5) from "Ready for handoff" to
# "In Translation" to reflect the newly generated report.
# ------------------------------------------------------------------

# ---------------- Sheet "Overview" ----------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A5").Value = "969e8835-09f2-4e4a-b193-79522ee0a29a.md"
$ws.Range("E5").Value = "In Translation"
$ws.Range("F5").Value = "In Translation"
$ws.Range("G5").Value = "2016-09-05 20:48:23"

$ws.Range("A6").Value = "8df1a12c-6253-4bfa-bf3e-77a8bf70f174.md"
$ws.Range("E6").Value = "Ready for handoff"
$ws.Range("F6").Value = "Ready for handoff"
$ws.Range("G6").Value = "2016-09-05 20:46:36"

# Rebuild the hyperlinks so B5/B6 keep referencing the same
# relationship ids (and therefore the same targets) while their
# display text is refreshed to match the new cell content.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f6f03daf86a8521f6b010318875c8c141fd481a3/e2e/ceeec683-33a2-4252-b233-b67ccb126543.md", [Type]::Missing, [Type]::Missing, "e2e\ceeec683-33a2-4252-b233-b67ccb126543.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/992e10ae58e476cf1c5b2595cda5342f2f3554f3/e2e/25b6cd32-4a59-414e-a9e1-a49af69ad4d9.md", [Type]::Missing, [Type]::Missing, "e2e\25b6cd32-4a59-414e-a9e1-a49af69ad4d9.md")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/992e10ae58e476cf1c5b2595cda5342f2f3554f3/e2e/63c391f1-7e21-4f53-8456-995ee4af4bd1.md", [Type]::Missing, [Type]::Missing, "e2e\63c391f1-7e21-4f53-8456-995ee4af4bd1.md")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b5d533a5e563407926f2cf142d2472b48657d1b/e2e/8df1a12c-6253-4bfa-bf3e-77a8bf70f174.md", [Type]::Missing, [Type]::Missing, "e2e\969e8835-09f2-4e4a-b193-79522ee0a29a.md")
$ws.Hyperlinks.Add($ws.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/40be608ee3f4575f747071b4b1d47dce1d671213/e2e/969e8835-09f2-4e4a-b193-79522ee0a29a.md", [Type]::Missing, [Type]::Missing, "e2e\8df1a12c-6253-4bfa-bf3e-77a8bf70f174.md")
$ws.Hyperlinks.Add($ws.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a2dad031c26c5cd43ec4096a6b9d97d717989830/e2e/bf4a510b-5d1c-4f1e-96ae-c1ad50c1ab36.md", [Type]::Missing, [Type]::Missing, "e2e\bf4a510b-5d1c-4f1e-96ae-c1ad50c1ab36.md")

# ---------------- Sheet "zh-cn" ----------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A5").Value = "969e8835-09f2-4e4a-b193-79522ee0a29a.md"
$ws.Range("C5").Value = "In Translation"
$ws.Range("G5").Value = "969e8835-09f2-4e4a-b193-79522ee0a29a.b8c306e62e59bbcd0f22bd7750ecceb5593be0b2.zh-cn.xlf"
$ws.Range("H5").Value = "2016-09-05 20:48:18"

$ws.Range("A6").Value = "8df1a12c-6253-4bfa-bf3e-77a8bf70f174.md"
$ws.Range("C6").Value = "Ready for handoff"
$ws.Range("G6").Value = "8df1a12c-6253-4bfa-bf3e-77a8bf70f174.7b9d39b001ce106c479d51f407d46885c0a499eb.zh-cn.xlf"
$ws.Range("H6").Value = "2016-09-05 20:46:31"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f6f03daf86a8521f6b010318875c8c141fd481a3/e2e/ceeec683-33a2-4252-b233-b67ccb126543.md", [Type]::Missing, [Type]::Missing, "ceeec683-33a2-4252-b233-b67ccb126543.md")
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/5cdc92c73d15a5da5e32ee194b981c8ec8a1f6a7/e2e/ceeec683-33a2-4252-b233-b67ccb126543.md", [Type]::Missing, [Type]::Missing, "ceeec683-33a2-4252-b233-b67ccb126543.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/992e10ae58e476cf1c5b2595cda5342f2f3554f3/e2e/25b6cd32-4a59-414e-a9e1-a49af69ad4d9.md", [Type]::Missing, [Type]::Missing, "25b6cd32-4a59-414e-a9e1-a49af69ad4d9.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/992e10ae58e476cf1c5b2595cda5342f2f3554f3/e2e/63c391f1-7e21-4f53-8456-995ee4af4bd1.md", [Type]::Missing, [Type]::Missing, "63c391f1-7e21-4f53-8456-995ee4af4bd1.md")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b5d533a5e563407926f2cf142d2472b48657d1b/e2e/8df1a12c-6253-4bfa-bf3e-77a8bf70f174.md", [Type]::Missing, [Type]::Missing, "969e8835-09f2-4e4a-b193-79522ee0a29a.md")
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/40be608ee3f4575f747071b4b1d47dce1d671213/e2e/969e8835-09f2-4e4a-b193-79522ee0a29a.md", [Type]::Missing, [Type]::Missing, "8df1a12c-6253-4bfa-bf3e-77a8bf70f174.md")
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a2dad031c26c5cd43ec4096a6b9d97d717989830/e2e/bf4a510b-5d1c-4f1e-96ae-c1ad50c1ab36.md", [Type]::Missing, [Type]::Missing, "bf4a510b-5d1c-4f1e-96ae-c1ad50c1ab36.md")

# ---------------- Sheet "de-de" ----------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A5").Value = "969e8835-09f2-4e4a-b193-79522ee0a29a.md"
$ws.Range("C5").Value = "In Translation"
$ws.Range("G5").Value = "969e8835-09f2-4e4a-b193-79522ee0a29a.b8c306e62e59bbcd0f22bd7750ecceb5593be0b2.de-de.xlf"
$ws.Range("H5").Value = "2016-09-05 20:48:23"

$ws.Range("A6").Value = "8df1a12c-6253-4bfa-bf3e-77a8bf70f174.md"
$ws.Range("C6").Value = "Ready for handoff"
$ws.Range("G6").Value = "8df1a12c-6253-4bfa-bf3e-77a8bf70f174.7b9d39b001ce106c479d51f407d46885c0a499eb.de-de.xlf"
$ws.Range("H6").Value = "2016-09-05 20:46:36"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f6f03daf86a8521f6b010318875c8c141fd481a3/e2e/ceeec683-33a2-4252-b233-b67ccb126543.md", [Type]::Missing, [Type]::Missing, "ceeec683-33a2-4252-b233-b67ccb126543.md")
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/e9cd0024d848c7c609f9b15c6856ce9d2dfa335f/e2e/ceeec683-33a2-4252-b233-b67ccb126543.md", [Type]::Missing, [Type]::Missing, "ceeec683-33a2-4252-b233-b67ccb126543.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/992e10ae58e476cf1c5b2595cda5342f2f3554f3/e2e/25b6cd32-4a59-414e-a9e1-a49af69ad4d9.md", [Type]::Missing, [Type]::Missing, "25b6cd32-4a59-414e-a9e1-a49af69ad4d9.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/992e10ae58e476cf1c5b2595cda5342f2f3554f3/e2e/63c391f1-7e21-4f53-8456-995ee4af4bd1.md", [Type]::Missing, [Type]::Missing, "63c391f1-7e21-4f53-8456-995ee4af4bd1.md")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b5d533a5e563407926f2cf142d2472b48657d1b/e2e/8df1a12c-6253-4bfa-bf3e-77a8bf70f174.md", [Type]::Missing, [Type]::Missing, "969e8835-09f2-4e4a-b193-79522ee0a29a.md")
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/40be608ee3f4575f747071b4b1d47dce1d671213/e2e/969e8835-09f2-4e4a-b193-79522ee0a29a.md", [Type]::Missing, [Type]::Missing, "8df1a12c-6253-4bfa-bf3e-77a8bf70f174.md")
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a2dad031c26c5cd43ec4096a6b9d97d717989830/e2e/bf4a510b-5d1c-4f1e-96ae-c1ad50c1ab36.md", [Type]::Missing, [Type]::Missing, "bf4a510b-5d1c-4f1e-96ae-c1ad50c1ab36.md")
